# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.066.80'
$ws.Range("E2").Value = '  +4.62%  '
$ws.Range("D3").Value = '2.468.79'
$ws.Range("E3").Value = '  +5.92%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''566.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.40%  '
$ws.Range("D6").Value = '''143.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.47%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +1.85%  '
$ws.Range("D9").Value = '2.468.23'
$ws.Range("E9").Value = '  +5.77%  '
$ws.Range("E10").Value = '  +4.15%  '
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("E13").Value = '  +4.81%  '
$ws.Range("D14").Value = '''26.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.90%  '
$ws.Range("D15").Value = '2.909.69'
$ws.Range("E15").Value = '  +6.02%  '
$ws.Range("D16").Value = '62.958.36'
$ws.Range("E16").Value = '  +4.50%  '
$ws.Range("E17").Value = '  +4.77%  '
$ws.Range("D18").Value = '2.467.55'
$ws.Range("E18").Value = '  +6.01%  '
$ws.Range("D19").Value = '''11.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.03%  '
$ws.Range("D20").Value = '''341.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.39%  '
$ws.Range("D21").Value = '''4.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.60%  '
$ws.Range("E22").Value = '  +3.24%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '''65.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.31%  '
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '''1.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.91%  '
$ws.Range("D28").Value = '''8.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("E29").Value = '  +6.01%  '
$ws.Range("E30").Value = '  +12.72%  '
$ws.Range("D31").Value = '0.0₃0805'
$ws.Range("E31").Value = '  +9.87%  '
$ws.Range("D33").Value = '''177.49'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.57%  '
$ws.Range("D34").Value = '''1.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.48%  '
$ws.Range("E35").Value = '  +3.84%  '
$ws.Range("D36").Value = '''18.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.79%  '
$ws.Range("D37").Value = '''368.93'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.88%  '
$ws.Range("D38").Value = '''4.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.55%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '''1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("E41").Value = '  +10.98%  '
$ws.Range("D42").Value = '''40.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.33%  '
$ws.Range("D43").Value = '''150.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.21%  '
$ws.Range("D44").Value = '''3.72'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.89%  '
$ws.Range("D45").Value = '''20.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.01%  '
$ws.Range("D46").Value = '''0.598'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.34%  '
$ws.Range("D47").Value = '''0.0959'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").Value = '''0.0517'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.89%  '
$ws.Range("D49").Value = '0.0₆0240'
$ws.Range("E49").Value = '  +9.60%  '
$ws.Range("D50").Value = '''0.0226'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.47%  '
$ws.Range("D51").Value = '''18.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.78%  '
